$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin (D2 has two dots, stays text naturally)
$ws.Range("D2").Value = "27.643.06"
$ws.Range("E2").Value = "  +0.06%  "

# Row 3 - Ethereum (D3 has two dots, stays text naturally)
$ws.Range("D3").Value = "1.755.55"
$ws.Range("E3").Value = "  -0.35%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB (plain decimal -> force text with leading apostrophe)
$ws.Range("D5").Value = "'324.36"
$ws.Range("E5").Value = "  -0.04%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.12%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.4492"
$ws.Range("E7").Value = "  +5.13%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3550"
$ws.Range("E8").Value = "  -1.61%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.07465"
$ws.Range("E9").Value = "  -1.42%  "

# Row 10 - OKB
$ws.Range("D10").Value = "'41.65"
$ws.Range("E10").Value = "  -1.56%  "

# Row 11 - Polygon
$ws.Range("D11").Value = "'1.089"
$ws.Range("E11").Value = "  -1.91%  "

# Row 12 - BinanceUSD
$ws.Range("E12").Value = "  +0.08%  "

# Row 13 - Solana
$ws.Range("D13").Value = "'20.74"
$ws.Range("E13").Value = "  -0.63%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'5.995"
$ws.Range("E14").Value = "  -1.41%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "'7.169"
$ws.Range("E15").Value = "  -0.82%  "

# Row 16 - WrappedEther (two dots, stays text naturally)
$ws.Range("D16").Value = "1.755.60"
$ws.Range("E16").Value = "  -0.26%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "'93.34"
$ws.Range("E17").Value = "  +0.05%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "'0.00001061"
$ws.Range("E18").Value = "  -0.91%  "

# Row 19 - TRON
$ws.Range("D19").Value = "'0.06474"
$ws.Range("E19").Value = "  +0.89%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.13%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "'17.07"
$ws.Range("E21").Value = "  -0.52%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'5.749"
$ws.Range("E22").Value = "  -2.64%  "

# Row 23 - WrappedBTC (two dots, stays text naturally)
$ws.Range("D23").Value = "27.683.54"
$ws.Range("E23").Value = "  +0.03%  "

# Row 24 - Cosmos
$ws.Range("D24").Value = "'11.24"
$ws.Range("E24").Value = "  -0.62%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "'2.110"
$ws.Range("E25").Value = "  -0.70%  "

# Row 26 - Monero
$ws.Range("D26").Value = "'164.13"
$ws.Range("E26").Value = "  +0.82%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'20.15"
$ws.Range("E27").Value = "  -1.12%  "

# Row 28 - WrappedliquidstakedEther2.0 (two dots, stays text naturally)
$ws.Range("D28").Value = "1.954.84"
$ws.Range("E28").Value = "  -0.17%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "'2.084"
$ws.Range("E29").Value = "  -3.84%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "'125.34"
$ws.Range("E30").Value = "  -0.29%  "

# Row 31 - ImmutableX
$ws.Range("D31").Value = "'1.071"
$ws.Range("E31").Value = "  -3.65%  "

# Row 32 - Stellar
$ws.Range("D32").Value = "'0.09181"
$ws.Range("E32").Value = "  +2.58%  "

# Row 33 - HuobiToken
$ws.Range("D33").Value = "'3.655"
$ws.Range("E33").Value = "  -0.06%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "'5.485"
$ws.Range("E34").Value = "  -2.35%  "

# Row 35 - VeChain
$ws.Range("E35").Value = "  -0.20%  "

# Row 36 - Aptos
$ws.Range("D36").Value = "'11.71"
$ws.Range("E36").Value = "  -4.49%  "

# Row 37 - Hedera
$ws.Range("D37").Value = "'0.06043"

# Row 38 - Algorand
$ws.Range("D38").Value = "'0.2075"
$ws.Range("E38").Value = "  -1.72%  "

# Row 39 - TheSandbox
$ws.Range("D39").Value = "'0.6297"
$ws.Range("E39").Value = "  -1.24%  "

# Row 40 - InternetComputer(DFINITY)
$ws.Range("D40").Value = "'4.950"
$ws.Range("E40").Value = "  -0.33%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "'1.181"
$ws.Range("E41").Value = "  -0.70%  "

# Row 42 - WEMIXTOKEN
$ws.Range("D42").Value = "'1.389"
$ws.Range("E42").Value = "  -0.46%  "

# Row 43 - FraxShare
$ws.Range("D43").Value = "'7.766"
$ws.Range("E43").Value = "  -1.84%  "

# Row 44 - EnergySwap
$ws.Range("D44").Value = "'13.15"
$ws.Range("E44").Value = "  -2.11%  "

# Row 45 - was Decentraland, becomes PancakeSwap
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "'3.710"
$ws.Range("E45").Value = "  -0.07%  "

# Row 46 - was PancakeSwap, becomes Decentraland
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5880"
$ws.Range("E46").Value = "  -1.23%  "

# Row 47 - Quant
$ws.Range("D47").Value = "'123.21"
$ws.Range("E47").Value = "  -0.24%  "

# Row 48 - NEARProtocol
$ws.Range("D48").Value = "'1.944"
$ws.Range("E48").Value = "  -2.54%  "

# Row 49 - Cronos
$ws.Range("D49").Value = "'0.06897"
$ws.Range("E49").Value = "  +0.46%  "

# Row 50 - EOS
$ws.Range("D50").Value = "'1.129"
$ws.Range("E50").Value = "  -3.84%  "

# Row 51 - Aave
$ws.Range("D51").Value = "'71.69"
$ws.Range("E51").Value = "  -2.46%  "
